$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0.005773782730102539"
$ws.Range("C2").Value = [double]"0.001387119334192364"
$ws.Range("D2").Value = [double]"0.002453947067260742"
$ws.Range("E2").Value = [double]"0.0009246126602285256"

$ws.Range("B3").Value = [double]"0.005405712127685547"
$ws.Range("C3").Value = [double]"0.001822770543089027"
$ws.Range("D3").Value = [double]"0.001738262176513672"
$ws.Range("E3").Value = [double]"0.001064065696672684"

$ws.Range("B4").Value = [double]"0.001746988296508789"
$ws.Range("C4").Value = [double]"0.0005380406310797586"
$ws.Range("D4").Value = [double]"0.0007751941680908203"
$ws.Range("E4").Value = [double]"5.536582506633495E-05"

$ws.Range("B5").Value = [double]"0.003672170639038086"
$ws.Range("C5").Value = [double]"0.00144938825540057"
$ws.Range("D5").Value = [double]"0.001158428192138672"
$ws.Range("E5").Value = [double]"0.0004784638648038072"

$ws.Range("B6").Value = [double]"0.00190424919128418"
$ws.Range("C6").Value = [double]"0.0004895047323137331"
$ws.Range("D6").Value = [double]"0.0008974552154541015"
$ws.Range("E6").Value = [double]"0.0003121654216632886"

$ws.Range("B7").Value = [double]"0.002164363861083984"
$ws.Range("C7").Value = [double]"6.449599351999547E-05"
$ws.Range("D7").Value = [double]"0.0007542133331298828"
$ws.Range("E7").Value = [double]"3.879257114415096E-05"

$ws.Range("B8").Value = [double]"0.001489782333374024"
$ws.Range("C8").Value = [double]"6.328270469304457E-05"
$ws.Range("D8").Value = [double]"0.0007410049438476562"
$ws.Range("E8").Value = [double]"9.298703601644034E-05"

$ws.Range("B9").Value = [double]"0.001821279525756836"
$ws.Range("C9").Value = [double]"0.000198921523534491"
$ws.Range("D9").Value = [double]"0.0006938934326171875"
$ws.Range("E9").Value = [double]"1.903983963114373E-05"

$ws.Range("B10").Value = [double]"0.00189824104309082"
$ws.Range("C10").Value = [double]"0.0001192670110194489"
$ws.Range("D10").Value = [double]"0.0007424354553222656"
$ws.Range("E10").Value = [double]"3.778313937817148E-05"

$ws.Range("B11").Value = [double]"0.001833057403564453"
$ws.Range("C11").Value = [double]"0.0001682022106085905"
$ws.Range("D11").Value = [double]"0.0007028579711914062"
$ws.Range("E11").Value = [double]"5.308099329067785E-05"

$ws.Range("B12").Value = [double]"0.001828527450561523"
$ws.Range("C12").Value = [double]"0.000630344918018385"
$ws.Range("D12").Value = [double]"0.0008051395416259766"
$ws.Range("E12").Value = [double]"0.0002328194231190658"

$ws.Range("B13").Value = [double]"0.001802396774291992"
$ws.Range("C13").Value = [double]"0.0002787191872836242"
$ws.Range("D13").Value = [double]"0.0007100582122802734"
$ws.Range("E13").Value = [double]"7.823318319715998E-05"

$ws.Range("B14").Value = [double]"0.00199127197265625"
$ws.Range("C14").Value = [double]"0.0005310252219676476"
$ws.Range("D14").Value = [double]"0.000689697265625"
$ws.Range("E14").Value = [double]"4.229908447011366E-05"

$ws.Range("B15").Value = [double]"0.001943016052246094"
$ws.Range("C15").Value = [double]"0.0004175408341310815"
$ws.Range("D15").Value = [double]"0.001000738143920898"
$ws.Range("E15").Value = [double]"0.0005815304819620885"

$ws.Range("B16").Value = [double]"0.00172266960144043"
$ws.Range("C16").Value = [double]"0.0002592130212230674"
$ws.Range("D16").Value = [double]"0.0006938457489013671"
$ws.Range("E16").Value = [double]"3.536749772973585E-05"

$ws.Range("B17").Value = [double]"0.001535272598266601"
$ws.Range("C17").Value = [double]"0.0001610624777795063"
$ws.Range("D17").Value = [double]"0.0006433486938476563"
$ws.Range("E17").Value = [double]"4.158216874248606E-05"

